# "new preparation file versuch 7.csv"
# Versuch_5 -> Versuch_6: rename the sheet, bold the header row, drop the
# trailing rows that belonged to the old (longer) data pull, and leave the
# view scrolled/selected where the author left off.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Sheet tab rename: Versuch_5 -> Versuch_6
$ws.Name = "Versuch_6"

# 2) Header row (A1:K1) becomes bold (new cell style fontId=16 in the diff)
$ws.Range("A1:K1").Font.Bold = $true

# 3) The new export is shorter: rows 378-386 (the tail of the old 386-row
#    dataset) are gone, shrinking the used range to A1:K377
$ws.Range("A378:K386").EntireRow.Delete()

# 4) Leave the view scrolled down a bit with the post-edit selection
$win = $excel.ActiveWindow
$win.ScrollRow = 4
$win.ScrollColumn = 1
$ws.Range("C390:D390").Select()
